# Update the "Кількість" column header to "Кількість упак." and move the
# active selection from A3 to G3, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Кількість упак."

[void]$ws.Range("G3").Select()
